$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting the existing D:K data to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy formats (number format, font, alignment) from column F (the old column D,
# now shifted right by two) onto the two freshly inserted, unformatted columns.
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E7:E102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$data = @{
  7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
  8 = @(1264700, 1267400, 2549200, 1253000, 1094200, 1077000, 2194000, 1091500, 1046300, 1053800)
  9 = @(1093500, 1087100, 2160500, 1065800, 946900, 939100, 1849300, 929100, 909900, 915700)
  10 = @(171300, 180300, 388700, 187200, 147300, 137900, 344700, 162400, 136400, 138100)
  11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
  13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  14 = @(3900, 0, 10300, 10300, 0, 0, 19700, 19700, 0, 10900)
  15 = @(52600, 50500, 98500, 46800, 40400, 38800, 80900, 42500, 37400, 37200)
  16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  17 = @(1180300, 1167600, 2330300, 1154600, 1017900, 1004900, 2006300, 1019500, 990500, 1008500)
  18 = @(84400, 99800, 218900, 98400, 76300, 72100, 187700, 72000, 55800, 45300)
  19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  20 = @(7000, 7600, 16300, 5100, 5500, 4400, 11100, 5500, 11000, 4200)
  21 = @(144000, 157900, 333800, 150200, 122200, 115300, 279700, 120100, 104200, 86700)
  22 = @(50500, 50700, 97300, 47200, 38500, 37700, 78500, 40900, 42400, 44500)
  23 = @(40900, 56700, 137900, 56300, 43300, 38800, 120300, 36700, 24400, 5100)
  24 = @(11200, 14100, 33400, 12300, -6300, 14000, 45600, 13200, 3900, 1100)
  25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  26 = @(29700, 42700, 104500, 44000, 49600, 24800, 74800, 23500, 20500, 4000)
  27 = @(23900, 31800, 77600, 32600, 26000, 17900, 56100, 15400, 19500, 6300)
  28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  29 = @("NA", "NA", "NA", "NA", 71500, "NA", "NA", "NA", "NA", "NA")
  30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  32 = @(-7000, -7600, -16300, -5100, -5500, -4400, -11100, -5500, -11000, -4200)
  33 = @(23900, 31800, 77600, 32600, 97500, 17900, 56100, 15400, 19500, 6300)
  34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  35 = @(23900, 31800, 77600, 32600, 97500, 17900, 56100, 15400, 19500, 6300)
  38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
  39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  41 = @(175200, 160400, 141000, 119700, 122500, 107300, 73800, 65200, 99000, 68200)
  42 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  43 = @(706700, 746800, 775600, 806400, 691700, 716400, 714200, 691500, 573800, 592700)
  44 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  45 = @(110700, 116100, 102700, 115300, 106500, 80300, 83200, 87500, 135300, 144800)
  46 = @(992500, 1023300, 1019300, 1041300, 920800, 904100, 871200, 844200, 762900, 805700)
  47 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  48 = @(979800, 965400, 965800, 973500, 912600, 946100, 911500, 897100, 1784400, 863500)
  49 = @(3758400, 3756300, 3766500, 3743300, 3109300, 3098900, 3102000, 3096800, 3091600, 3012800)
  50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  52 = @(233500, 235200, 213100, 210600, 184400, 174800, 169400, 164700, 173900, 163300)
  53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  54 = @(5964300, 5980200, 5964800, 5968600, 5127200, 5123800, 5054200, 5002900, 4920600, 4845400)
  55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  57 = @(146700, 137500, 131800, 138400, 128200, 133300, 123600, 125100, 126600, 114200)
  58 = @(68900, 47500, 47800, 44000, 51700, 56500, 60700, 44300, 53000, 32800)
  59 = @(489500, 500100, 447700, 443200, 425600, 409700, 384500, 369100, 392100, 403300)
  60 = @(705200, 685100, 627300, 625700, 605400, 599500, 568900, 538500, 571600, 550300)
  61 = @(3249500, 3305500, 3386200, 3478000, 2677700, 2752700, 2734100, 2771400, 2685300, 2642100)
  62 = @(312800, 322400, 323100, 292100, 270600, 329600, 338700, 337900, 335600, 346500)
  63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  66 = @(5161200, 5205600, 5069300, 5116000, 4303800, 4402200, 4205900, 4201200, 4104900, 3876400)
  67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  72 = @(320400, 296600, 420500, 383600, 359700, 262500, 393300, 351200, 371700, 528600)
  73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  76 = @(803000, 774600, 895500, 852600, 823400, 721600, 848300, 801700, 815700, 969000)
  77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
  81 = @(23900, 31800, 77600, 32600, 97500, 17900, 56100, 15400, 19500, 6300)
  82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  83 = @(52600, 50500, 98500, 46800, 40400, 38800, 80900, 42500, 37000, 37300)
  84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  89 = @(113200, 164000, 217000, 50700, 108200, 89600, 40300, -55900, 65800, 102500)
  90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  91 = @(-46200, -39400, -81600, -39600, -59400, -68500, -105300, -50700, 74900, -76000)
  92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  94 = @(-50600, -50600, -596000, -556000, -23000, -70900, -99100, -41200, -91300, -31000)
  95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  96 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  100 = @(-47900, -94100, 397500, 502400, -69900, 14700, 33600, 63300, 55800, -81500)
  101 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  102 = @(14800, 19400, 18500, -2900, 15200, 33500, -25200, -33800, 30800, -10200)
}

# Write the updated / restated quarterly figures (columns D through M) for every
# data row on the sheet.
foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item([int]$r, 4 + $i).Value = $v
        }
    }
}

Write-Output "done"
